$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-11 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-12 Tuesday", 2)
$d.Content.Find.Execute("425÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "452÷4=", 2)
$d.Content.Find.Execute("843÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "627÷2=", 2)
$d.Content.Find.Execute("520÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "433÷2=", 2)
$d.Content.Find.Execute("224÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "236÷5=", 2)
$d.Content.Find.Execute("937÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "458÷3=", 2)
$d.Content.Find.Execute("895÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "573÷9=", 2)
$d.Content.Find.Execute("533÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "797÷2=", 2)
$d.Content.Find.Execute("383÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "586÷6=", 2)
$d.Content.Find.Execute("238÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "707÷8=", 2)
$d.Content.Find.Execute("604÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "478÷6=", 2)
$d.Content.Find.Execute("579÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "210÷2=", 2)
$d.Content.Find.Execute("110÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "429÷3=", 2)
$d.Content.Find.Execute("413÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "232÷5=", 2)
$d.Content.Find.Execute("935÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "201÷9=", 2)
$d.Content.Find.Execute("342÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "795÷4=", 2)
$d.Content.Find.Execute("220÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "280÷7=", 2)
$d.Content.Find.Execute("453÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "961÷6=", 2)
$d.Content.Find.Execute("146÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "723÷6=", 2)
$d.Content.Find.Execute("345÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "795÷7=", 2)
$d.Content.Find.Execute("538÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "332÷9=", 2)
$d.Content.Find.Execute("114÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "771÷7=", 2)
$d.Content.Find.Execute("793÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "152÷3=", 2)
$d.Content.Find.Execute("938÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "789÷6=", 2)
$d.Content.Find.Execute("535÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "279÷4=", 2)
$d.Content.Find.Execute("677÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "625÷7=", 2)
